$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source table previously began with 2009; that oldest year is dropped and every
# remaining year's row shifts up by one (2010 -> row 2, 2011 -> row 3, ... 2020 -> row 12).
$ws.Rows.Item(2).Delete()

# Seed row 13 from row 12 (copies formatting, including the bold/bordered year-cell
# style, and keeps column E's long-standing blank-cell representation) before
# overwriting it with the new 2021 figures.
$ws.Range("A12:U12").Copy($ws.Range("A13:U13"))

# Append the new 2021 data as row 13.
$ws.Cells.Item(13, 1).Value = "2021年"
$ws.Cells.Item(13, 2).Value = 62411
$ws.Cells.Item(13, 3).Value = 46817
$ws.Cells.Item(13, 4).Value = 114618
$ws.Cells.Item(13, 6).Value = 41442
$ws.Cells.Item(13, 7).Value = 63946
$ws.Cells.Item(13, 8).Value = 67750
$ws.Cells.Item(13, 9).Value = 62884
$ws.Cells.Item(13, 10).Value = 47193
$ws.Cells.Item(13, 11).Value = 60430
$ws.Cells.Item(13, 12).Value = 58288
$ws.Cells.Item(13, 13).Value = 58071
$ws.Cells.Item(13, 14).Value = 52579
$ws.Cells.Item(13, 15).Value = 56171
$ws.Cells.Item(13, 16).Value = 43366
$ws.Cells.Item(13, 17).Value = 59271
$ws.Cells.Item(13, 18).Value = 77708
$ws.Cells.Item(13, 19).Value = 64490
$ws.Cells.Item(13, 20).Value = 62665
$ws.Cells.Item(13, 21).Value = 95416
